{"js": "// Apply the \"Added many more features\" edits to the Diamond Queen review.\n// Each entry is an exact, unique (or intentionally-repeated) text replacement.\nconst replacements = [\n  {\n    find: \"Play Diamond Queen Free Slot by IGT\",\n    replace: \"Play Diamond Queen Slot Free\",\n  },\n  {\n    find: \"Well-designed and elegant graphics\",\n    replace: \"Fantasy-themed slot with elegant design\",\n  },\n  {\n    find: \"Stacked symbols feature for increased winning potential\",\n    replace: \"Special features like free spins and stacked symbols\",\n  },\n  {\n    find: \"Mystical Diamond Bonus with free spins and extra Wilds\",\n    replace: \"Magical atmosphere with diamond and enchantment theme\",\n  },\n  {\n    find: \"Medium volatility with a high RTP of 96.08%\",\n    replace: \"Medium volatility with a decent RTP\",\n  },\n  {\n    find: \"No background music, only short jingles for animations\",\n    replace: \"Lack of background music\",\n  },\n  {\n    find: \"The theme may not be appealing to players who do not enjoy fantasy slots\",\n    replace: \"Limited options for players who don't enjoy diamond-themed slots\",\n  },\n  {\n    find:\n      \"Read our review of Diamond Queen, a magical and elegant online slot by IGT. Play for free and trigger the Mystical Diamond Bonus for extra Wilds and free spins.\",\n    replace:\n      \"Read our review of Diamond Queen, a fantasy-themed slot with magical features. Play for free!\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items,text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Diamond Queen review.\n$d = $word.ActiveDocument\n\n# Disable smart-quote autocorrect so straight apostrophes in the new text\n# (e.g. \"don't\") are not turned into curly quotes on insertion.\n$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false\n$word.Options.AutoFormatReplaceQuotes = $false\n\n$replacements = @(\n  @{ Find = \"Play Diamond Queen Free Slot by IGT\"; Replace = \"Play Diamond Queen Slot Free\" },\n  @{ Find = \"Well-designed and elegant graphics\"; Replace = \"Fantasy-themed slot with elegant design\" },\n  @{ Find = \"Stacked symbols feature for increased winning potential\"; Replace = \"Special features like free spins and stacked symbols\" },\n  @{ Find = \"Mystical Diamond Bonus with free spins and extra Wilds\"; Replace = \"Magical atmosphere with diamond and enchantment theme\" },\n  @{ Find = \"Medium volatility with a high RTP of 96.08%\"; Replace = \"Medium volatility with a decent RTP\" },\n  @{ Find = \"No background music, only short jingles for animations\"; Replace = \"Lack of background music\" },\n  @{ Find = \"The theme may not be appealing to players who do not enjoy fantasy slots\"; Replace = \"Limited options for players who don't enjoy diamond-themed slots\" },\n  @{ Find = \"Play Diamond Queen Free Slot by IGT\"; Replace = \"Play Diamond Queen Slot Free\" },\n  @{ Find = \"Read our review of Diamond Queen, a magical and elegant online slot by IGT. Play for free and trigger the Mystical Diamond Bonus for extra Wilds and free spins.\"; Replace = \"Read our review of Diamond Queen, a fantasy-themed slot with magical features. Play for free!\" }\n)\n\n# wdFindStop = 0\n$wdFindStop = 0\n\nforeach ($pair in $replacements) {\n  $searchRange = $d.Content\n  $find = $searchRange.Find\n  $find.ClearFormatting()\n  $find.Text = $pair.Find\n  $find.Forward = $true\n  $find.Wrap = $wdFindStop\n  $find.MatchCase = $true\n  $found = $find.Execute()\n  if ($found) {\n    # Assigning Range.Text directly (rather than Find's own Replacement/\n    # Execute-with-replace) rewrites only the matched text and keeps a\n    # straight apostrophe instead of Find's smart-quote substitution.\n    $searchRange.Text = $pair.Replace\n  }\n}\n"}
